$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBSSports")

$values = @(
  "Menu",
  "NFL",
  "HOME",
  "SCORES",
  "SCHEDULE",
  "STANDINGS",
  "TEAMS",
  "STATS",
  "PLAYERS",
  "PLAY",
  "WATCH",
  "BET",
  "POSCASTS",
  "LOG IN",
  "Dont have an Account?Sign Up?",
  "REGISTER",
  "EMAIL",
  "PASSWORD",
  "CONFIRM PASSWORD",
  "FIRST NAME",
  "LAST NAME",
  "MONTH",
  "DAY",
  "YEAR",
  "ZIP CODE",
  "CHECK BOX",
  "REGISTER",
  "REGISTER"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Columns.Item(1).ColumnWidth = 24.5

$ws.Activate()
$ws.Range("A14").Select()
